# Logged Week 17 data and fixed Simulate_Season.py tiebreaking method
$wb = $excel.ActiveWorkbook

# --- OFF sheet: update row 3 (R) ---
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 246
$wsOff.Range("C3").Value = 167
$wsOff.Range("D3").Value = 57
$wsOff.Range("E3").Value = 26
$wsOff.Range("F3").Value = 10
$wsOff.Range("G3").Value = 2

# --- DEF sheet: update row 3 (R) ---
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 223
$wsDef.Range("C3").Value = 150
$wsDef.Range("D3").Value = 60
$wsDef.Range("E3").Value = 32
